# "First set of edits after R&R"
#
# This script updates the SS_att.xlsx workbook:
#  - Attrition sheet: drops the old "Commitment arms" merged banner row,
#    turning the former two-row header into a single header row, adds a
#    "Structure" column header (replacing the removed "Forced" label),
#    and tightens the precision of the Take-up ROUND() formulas from 2 to
#    3 decimal places.
#  - SS_att sheet: refreshes the underlying summary-stats numbers (Take-up
#    row) with the values from the revised analysis.

$wb = $excel.ActiveWorkbook

$wsAtt = $wb.Worksheets.Item("Attrition")
$wsSS  = $wb.Worksheets.Item("SS_att")

# ------------------------------------------------------------------
# 1) Update the underlying stats sheet (SS_att) with the revised numbers
# ------------------------------------------------------------------
$wsSS.Range("B2").Value = 0.96716697936210128
$wsSS.Range("C2").Value = 0.95470085470085497
$wsSS.Range("D2").Value = 0.96133682830930534
$wsSS.Range("E2").Value = 0.96092503987240829
$wsSS.Range("L2").Value = 0.82203658081697284

$wsSS.Range("B3").Value = 0.010616477679248965
$wsSS.Range("C3").Value = 0.014476326826305447
$wsSS.Range("D3").Value = 0.013703690061123075
$wsSS.Range("E3").Value = 0.009593410319644715

$wsSS.Range("B4").Value = 2635
$wsSS.Range("C4").Value = 2535
$wsSS.Range("D4").Value = 3494
$wsSS.Range("E4").Value = 8664

# Cosmetic: selection left on O2 in the SS_att sheet after the edits
$wsSS.Range("O2").Select()

# ------------------------------------------------------------------
# 2) Remove the old "Commitment arms" banner row above the header on the
#    Attrition sheet - this merges the two header rows into one and
#    shifts everything below up by one row.
# ------------------------------------------------------------------
$wsAtt.Rows("2").Delete()

# ------------------------------------------------------------------
# 3) Rebuild the (now single) header row with the new "Structure" column
# ------------------------------------------------------------------
$wsAtt.Range("B2").Value = "Control"
$wsAtt.Range("C2").Value = "Structure"
$wsAtt.Range("D2").Value = "Choice"
$wsAtt.Range("E2").Value = "p-value"

# ------------------------------------------------------------------
# 4) Tighten the Take-up row formulas to 3 decimal places
# ------------------------------------------------------------------
$wsAtt.Range("B3").Formula = "=ROUND(SS_att!B2,3)"
$wsAtt.Range("C3").Formula = "=ROUND(SS_att!C2,3)"
$wsAtt.Range("D3").Formula = "=ROUND(SS_att!D2,3)"

# Cosmetic: selection on the Attrition sheet now spans one fewer row
$wsAtt.Range("A2:E20").Select()

$wb.Save()
